$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 367, pushing the existing row 367 (and
# everything below it) down by one.
$ws.Rows.Item(367).Insert()

# Populate the newly inserted row 367 with the new data record.
$ws.Range("A367").Value = 10
$ws.Range("B367").Value = "Vega Modelo de Temuco"
$ws.Range("C367").Value = "La Araucanía"
$ws.Range("D367").Value = 44706
$ws.Range("E367").Value = 9
$ws.Range("F367").Value = 100112008
$ws.Range("G367").Value = "Coliflor"
$ws.Range("H367").Value = "Sin especificar"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 3000
$ws.Range("K367").Value = 1200
$ws.Range("L367").Value = 1200
$ws.Range("M367").Value = 1200
$ws.Range("N367").Value = "$/unidad"
$ws.Range("O367").Value = "Región Metropolitana"
$ws.Range("P367").Value = 1200
$ws.Range("Q367").Value = 1
$ws.Range("R367").Value = "Hortaliza"
